# Update the "Date & Time Report Created" label to clarify the timestamp
# is in UTC.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Threshold Report")

$ws.Range("E1").Value = "Date & Time Report Created (UTC)"
